# Fix multiple placeholders in one cell
#
# Translates the header/label cells from Russian to English, adds a new
# "Complex Content" column (D) that demonstrates a cell containing several
# rich-text runs / placeholders, and fixes the stray style on J3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: header labels, translated -------------------------------------
$ws.Range("A2").Value = "Title"
$ws.Range("C2").Value = "Conent"
$ws.Range("D2").Value = "Complex Content"
$ws.Range("F2").Value = "Date"
$ws.Range("H2").Value = "Author"

# --- Row 3: D3 - new rich-text placeholder cell ----------------------------
# Three runs: a red "From {author} in {date}", a plain separator, and an
# underlined "{foo}, {foo} and {bar}" - this is the "multiple placeholders
# in one cell" the commit message refers to.
$run1Text = "From {author} in {date}"
$run2Text = " `n"
$run3Text = "{foo}, {foo} and {bar}"
$fullText = $run1Text + $run2Text + $run3Text

$ws.Range("D3").Value = $fullText

$run1 = $ws.Range("D3").Characters(1, $run1Text.Length)
$run1.Font.Name = "Arial"
$run1.Font.Size = 10
$run1.Font.Color = 1974729

$run2 = $ws.Range("D3").Characters($run1Text.Length + 1, $run2Text.Length)
$run2.Font.Name = "Arial"
$run2.Font.Size = 10

$run3 = $ws.Range("D3").Characters($run1Text.Length + $run2Text.Length + 1, $run3Text.Length)
$run3.Font.Name = "Arial"
$run3.Font.Size = 10
$run3.Font.Underline = $true

# Cell-level formatting for D3: red font, wrap text, solid fill.
$ws.Range("D3").Font.Name = "Arial"
$ws.Range("D3").Font.Size = 10
$ws.Range("D3").Font.Color = 1974729
$ws.Range("D3").WrapText = $true
$ws.Range("D3").Interior.Color = 13093555

# J3 keeps the same text, just gets re-written.
$ws.Range("J3").Value = "{content}"

# --- Row 4: D4 - new, empty, plainly formatted cell ------------------------
$ws.Range("D4").Value = ""

# --- Column D sizing --------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 28

# --- Selection ---------------------------------------------------------------
$ws.Range("J3").Select()

Write-Output "done"
